$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_12_3_9", 0.6458022902904654, -23.02505197487528, 0.3509511836454897, -1.278055290244684, -0.2295123209787617, 0.2102670610279525, 14.26230868049242, 0.2401337672282458, 0.4683857702226017, 0.3542597687254238, 0.2697708819353354, 0.4585488643841051, -0.06259312912860371, 0.4780702388679436, 35.11875367406773, 54.62076687195893),
    @("model_12_3_12", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_22", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_21", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_20", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_19", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_18", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_17", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_16", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_15", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_14", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_13", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_24", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_23", 0.6383931775788736, -23.10748969621655, 0.3313169117313355, -1.388364236985702, -0.2814759963235094, 0.2146654303905574, 14.31124727300474, 0.2473980154060813, 0.4910661420305162, 0.3692320787182988, 0.2718993657310753, 0.4633200086231518, -0.08482046726337922, 0.4830445005948765, 35.07734920122535, 54.57936239911656),
    @("model_12_3_10", 0.6419327987711958, -23.19731153925137, 0.3389986513804766, -1.334093485931358, -0.2571804611668036, 0.2125641583471211, 14.3645693991318, 0.2445559409205664, 0.4799076562632908, 0.3622317985919286, 0.2710286096946006, 0.4610468071108627, -0.0742016036864126, 0.4806745241880335, 35.09702282910663, 54.59903602699784),
    @("model_12_3_11", 0.6394424168438981, -23.20432568632141, 0.3331906864870069, -1.37186992618305, -0.2743878612457455, 0.2140425566380626, 14.36873329982794, 0.2467047600149941, 0.4876747842778423, 0.3671897721464182, 0.2717672078397137, 0.462647335059938, -0.0816727494683056, 0.4823431899254507, 35.08316084187237, 54.58517403976358),
    @("model_12_3_8", 0.6466613230193305, -23.21738108578707, 0.3535338336325647, -1.266464750513067, -0.2237187130112821, 0.2097571020918164, 14.37648354891465, 0.2391782436139949, 0.4660026656936043, 0.3525904546537996, 0.2698003754854419, 0.4579924694706414, -0.06001603094200836, 0.4774901570711175, 35.12361014882786, 54.62562334671907),
    @("model_12_3_7", 0.6497370186875215, -23.70073748158287, 0.3627719291311196, -1.223673472362222, -0.2025197818381586, 0.2079312362800435, 14.66342478536814, 0.2357603517417901, 0.4572044482573012, 0.3464823999995457, 0.2692863323372691, 0.455994776592938, -0.05078894393743538, 0.4754074182718159, 35.14109569749363, 54.64310889538483),
    @("model_12_3_6", 0.6581278343715251, -24.5081548383484, 0.3885012404406331, -1.098186483378788, -0.1412274227722186, 0.2029500856256546, 15.14274260695787, 0.2262410732264324, 0.4314033536835134, 0.3288222134549729, 0.2673334797708395, 0.4504998175645076, -0.02561649688542467, 0.4696785274614119, 35.18959042758429, 54.6916036254755),
    @("model_12_3_5", 0.6736368436614425, -25.83206991164971, 0.4400746242175776, -0.8048830080860336, -0.00346601120452017, 0.1937432677568423, 15.92867578462285, 0.2071600570948178, 0.3710978928054295, 0.2891289749501237, 0.2657632140910283, 0.4401627741606988, 0.02091053098432749, 0.4589014147192732, 35.28244271682767, 54.78445591471888),
    @("model_12_3_4", 0.6748681855409469, -27.22608237780868, 0.4444539085274858, -0.7835597643828569, 0.006953688082644405, 0.193012290025969, 16.75622180273784, 0.2055398183506685, 0.3667136691352113, 0.2861267437429399, 0.2655559814887576, 0.4393316401375719, 0.02460455662284056, 0.4580348976000971, 35.29000282642727, 54.79201602431848),
    @("model_12_3_3", 0.6765810469609213, -28.76707086765381, 0.4509606861318993, -0.7436113646171385, 0.02538469592086878, 0.1919954614953115, 17.67101913754657, 0.2031324539440277, 0.3584999694618304, 0.280816211702929, 0.265369488407859, 0.438172867137288, 0.02974314088276375, 0.4568267932342004, 35.30056709055798, 54.80258028844919),
    @("model_12_3_2", 0.6785021157788209, -30.32404180327222, 0.4586743093999733, -0.6920248337996866, 0.04874298122517462, 0.1908550320591543, 18.59530434257203, 0.2002785832217208, 0.3478933801162974, 0.2740859816690091, 0.2651738177774576, 0.4368695824375443, 0.03550634733646274, 0.4554680250066173, 35.31248226714463, 54.81449546503584),
    @("model_12_3_1", 0.6804142929032229, -31.84404152097266, 0.4669713966866665, -0.6318623104702326, 0.07553575273398694, 0.1897198811163586, 19.49764183556792, 0.1972088436628908, 0.335523500443546, 0.2663661720532184, 0.2649572911351829, 0.4355684574396527, 0.04124287870966858, 0.4541115084238654, 35.32441320919956, 54.82642640709076),
    @("model_12_3_0", 0.6625621858825883, -44.13982607852679, 0.4781908929206365, 0.08105672875285574, 0.3371059926404703, 0.2003176630146764, 26.79695069916838, 0.1930578771574766, 0.1889418372491477, 0.1909998572033121, 0.2877405289119224, 0.4475686126335004, -0.01231344235223508, 0.4666225350680611, 35.21570171479812, 54.71771491268932)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value = $rowVals[$j]
    }
}